$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Valid/Invalid results for rows 2 and 3
$ws.Range("C2").Value = "Invalid"
$ws.Range("C3").Value = "Valid"

# Clear out the old test-data row (laksh@yahoo.com / Lakshmi), keeping the
# existing cell formatting in place (A4/B4 keep their style, just no value).
$ws.Range("A4").Value = ""
$ws.Range("B4").Value = ""

# The hyperlink engine can only clear ALL hyperlinks on a sheet at once, so
# remove them all and re-create the ones that should remain (everything
# except the A4 link, which belonged to the row we just cleared).
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:pavanoltraining@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:lakshmi@yahoo.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:laks@yahoo.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:test@123")
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:abc123@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:test@123")
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:bevafi7510@aravites.com")
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:Sowmi@96")

# Re-adding hyperlinks resets the cell font to the default hyperlink font;
# restore the workbook's original 16pt font size so the cells reuse the
# original style entries instead of creating new ones.
$ws.Range("A2").Font.Size = 16
$ws.Range("A3").Font.Size = 16
$ws.Range("A5").Font.Size = 16
$ws.Range("B2").Font.Size = 16
$ws.Range("A6").Font.Size = 16
$ws.Range("B6").Font.Size = 16
$ws.Range("A7").Font.Size = 16
$ws.Range("B7").Font.Size = 16

# Move the active selection to A4, matching the saved view state.
$ws.Range("A4").Select() | Out-Null
